# Apply scheduled-runner market data & profit updates across all class sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2289
$ws.Range("I96").Value = 1008
$ws.Range("J96").Value = 3249.75
$ws.Range("K96").Value = 3024
$ws.Range("L96").Value = 9749.25
$ws.Range("M96").Value = -1651
$ws.Range("N96").Value = -12495.25
$ws.Range("H113").Value = 4668.3335
$ws.Range("I113").Value = 4002.5
$ws.Range("K113").Value = 4002.5
$ws.Range("M113").Value = -748.5
$ws.Range("H129").Value = 2779488.2
$ws.Range("J129").Value = 1817.5758
$ws.Range("L129").Value = 5452.7274
$ws.Range("N129").Value = -15452.7274
$ws.Range("H131").Value = 2925
$ws.Range("I131").Value = 566.6667
$ws.Range("K131").Value = 1700.0001
$ws.Range("M131").Value = 3339.9999
$ws.Range("H137").Value = 8317.333000000001
$ws.Range("I137").Value = 8317.333000000001
$ws.Range("K137").Value = 24951.999
$ws.Range("M137").Value = -22401.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 990
$ws.Range("I2").Value = 990
$ws.Range("K2").Value = 990
$ws.Range("M2").Value = -877
$ws.Range("H17").Value = 999.5
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1846
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H110").Value = 497.5
$ws.Range("I110").Value = 497.5
$ws.Range("K110").Value = 497.5
$ws.Range("M110").Value = 1547.5
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 990
$ws.Range("I116").Value = 990
$ws.Range("K116").Value = 990
$ws.Range("M116").Value = 1304
$ws.Range("H120").Value = 69999
$ws.Range("J120").Value = 69999
$ws.Range("L120").Value = 69999
$ws.Range("N120").Value = -79675

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 990
$ws.Range("I3").Value = 990
$ws.Range("K3").Value = 990
$ws.Range("M3").Value = -876
$ws.Range("H92").Value = 19449.5
$ws.Range("J92").Value = 19449.5
$ws.Range("L92").Value = 19449.5
$ws.Range("N92").Value = -24441.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9939.817999999999
$ws.Range("I31").Value = 3867.8
$ws.Range("J31").Value = 14999.833
$ws.Range("K31").Value = 3867.8
$ws.Range("L31").Value = 14999.833
$ws.Range("M31").Value = -3572.8
$ws.Range("N31").Value = -15589.833
$ws.Range("H34").Value = 9939.817999999999
$ws.Range("I34").Value = 3867.8
$ws.Range("J34").Value = 14999.833
$ws.Range("K34").Value = 3867.8
$ws.Range("L34").Value = 14999.833
$ws.Range("M34").Value = -3665.8
$ws.Range("N34").Value = -15403.833
$ws.Range("H122").Value = 3999.5
$ws.Range("I122").Value = 3999.5
$ws.Range("K122").Value = 11998.5
$ws.Range("M122").Value = -9548.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 499
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 499
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1497
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -1945
$ws.Range("H13").Value = 150
$ws.Range("I13").Value = 150
$ws.Range("K13").Value = 450
$ws.Range("M13").Value = -282
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 4000
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H131").Value = 3497.6667
$ws.Range("J131").Value = 3497.6667
$ws.Range("L131").Value = 10493.0001
$ws.Range("N131").Value = -20573.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3374.75
$ws.Range("I70").Value = 3166.3333
$ws.Range("K70").Value = 3166.3333
$ws.Range("M70").Value = -2896.3333
$ws.Range("H73").Value = 3374.75
$ws.Range("I73").Value = 3166.3333
$ws.Range("K73").Value = 3166.3333
$ws.Range("M73").Value = -2230.3333
$ws.Range("H113").Value = 4337
$ws.Range("I113").Value = 4337
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4337
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2167
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1999
$ws.Range("I61").Value = 1999
$ws.Range("K61").Value = 1999
$ws.Range("M61").Value = -1797
$ws.Range("H94").Value = 23749.666
$ws.Range("J94").Value = 23749.666
$ws.Range("L94").Value = 23749.666
$ws.Range("N94").Value = -25101.666
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("M113").Value = 171
$ws.Range("H116").Value = 200000
$ws.Range("J116").Value = 200000
$ws.Range("L116").Value = 200000
$ws.Range("N116").Value = -209178
$ws.Range("H122").Value = 3989.25
$ws.Range("I122").Value = 3985.6667
$ws.Range("K122").Value = 11957.0001
$ws.Range("M122").Value = -9507.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 7924300
$ws.Range("I5").Value = 8000
$ws.Range("J5").Value = 11882450
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 11882450
$ws.Range("M5").Value = -7888
$ws.Range("N5").Value = -11882674
$ws.Range("H122").Value = 1600
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350
